$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped from
# 2023-09-21 (45190) to 2023-09-23 (45192) for every data row (2..218).
$ws.Range("C2:C218").Value = 45192
